# Rebrand the "Artificial Intelligence and Machine Learning" budget template
# to a "Product Development" budget template.
#
# The rename follows the (naive, non-word-boundary) substring-replacement
# pattern used by the original authoring tool:
#   "Artificial Intelligence and Machine Learning" -> "Product Development"
#   "ARTIFICIAL INTELLIGENCE AND MACHINE LEARNING" -> "PRODUCT DEVELOPMENT"
#   "AI/ML" -> "Product Development"
#   "AI"    -> "Product"
#   "ML"    -> "Product"
# (applied in that order), which is why e.g. "TOTAL TRAINING" becomes
# "TOTAL TRProductNING" (the "AI" inside "TRAINING" gets replaced too).

$wb = $excel.ActiveWorkbook

function Rebrand-Text($text) {
    $t = $text
    $t = $t.Replace("Artificial Intelligence and Machine Learning", "Product Development")
    $t = $t.Replace("ARTIFICIAL INTELLIGENCE AND MACHINE LEARNING", "PRODUCT DEVELOPMENT")
    $t = $t.Replace("AI/ML", "Product Development")
    $t = $t.Replace("AI", "Product")
    $t = $t.Replace("ML", "Product")
    return $t
}

function Set-Text($sheetName, $addr) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cell = $ws.Range($addr)
    $cell.Value2 = Rebrand-Text($cell.Value2)
}

# Touch a row with a value-preserving no-op so the worksheet XML grows an
# explicit (attribute-less, cell-less) <row r="N"/> element for currently
# fully-empty rows -- mirrors rows Excel re-serializes once the sheet is
# round-tripped through the app, even where no cell in that row has data.
function Touch-EmptyRow($sheetName, $rowNum) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = $ws.Rows.Item($rowNum)
    $row.OutlineLevel = 0
}

# ---------------------------------------------------------------------
# 1. Instructions & User Guide
# ---------------------------------------------------------------------
Set-Text "Instructions & User Guide" "A1"
Set-Text "Instructions & User Guide" "A56"
Set-Text "Instructions & User Guide" "B59"
foreach ($r in 2, 10, 20, 28, 37, 45, 54, 55, 60) {
    Touch-EmptyRow "Instructions & User Guide" $r
}

# ---------------------------------------------------------------------
# 2. Budget Summary
# ---------------------------------------------------------------------
Set-Text "Budget Summary" "A1"
foreach ($r in 2, 6) {
    Touch-EmptyRow "Budget Summary" $r
}

# ---------------------------------------------------------------------
# 3. Resources
# ---------------------------------------------------------------------
Set-Text "Resources" "A1"
Set-Text "Resources" "A5"
Set-Text "Resources" "A6"
foreach ($r in 2, 11) {
    Touch-EmptyRow "Resources" $r
}

# ---------------------------------------------------------------------
# 4. Logistics
# ---------------------------------------------------------------------
Set-Text "Logistics" "A1"
foreach ($r in 2, 9) {
    Touch-EmptyRow "Logistics" $r
}

# ---------------------------------------------------------------------
# 5. Technology
# ---------------------------------------------------------------------
Set-Text "Technology" "A1"
Set-Text "Technology" "A5"
foreach ($r in 2, 10) {
    Touch-EmptyRow "Technology" $r
}

# ---------------------------------------------------------------------
# 6. Training
# ---------------------------------------------------------------------
Set-Text "Training" "A1"
Set-Text "Training" "A4"
Set-Text "Training" "A10"
foreach ($r in 2, 9) {
    Touch-EmptyRow "Training" $r
}

# ---------------------------------------------------------------------
# 7. Contingency
# ---------------------------------------------------------------------
Set-Text "Contingency" "A1"
foreach ($r in 2, 5, 11, 13) {
    Touch-EmptyRow "Contingency" $r
}

# ---------------------------------------------------------------------
# 8. Timeline
# ---------------------------------------------------------------------
Set-Text "Timeline" "A1"
foreach ($r in 2) {
    Touch-EmptyRow "Timeline" $r
}
